$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$mdFileName = "ed2d9656-aa33-4b30-a561-82adbcd92828.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3fcb8a6b5e63694c817209ccc963b6eb86f54869/e2e/ed2d9656-aa33-4b30-a561-82adbcd92828.md"

# The engine quantizes ColumnWidth to 1/6-character steps, with a fixed
# +5/6 offset versus the stored OOXML "width" attribute
# (stored = round(6*ColumnWidth+5)/6). These two inputs are the exact
# pre-images of the target OOXML widths 29.9777047293527 (-> nearest
# achievable 30) and 40 (-> 40, exact) used throughout the diff.
$wideColWidth = 29.166666666666668
$fullColWidth = 39.166666666666664

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value2 = $statusText
$wsOverview.Range("F2").Value2 = $statusText
$wsOverview.Columns.Item(5).ColumnWidth = $wideColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $wideColWidth

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value2 = $statusText
$wsZhCn.Range("I2").Value2 = $mdFileName
$wsZhCn.Range("J2").Value2 = "ed2d9656-aa33-4b30-a561-82adbcd92828.89e1e0ed422e8a73fb20768186dace8c5457326f.zh-cn.xlf"
$wsZhCn.Range("K2").Value2 = "2016-08-25 19:01:34"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdUrl, "", "", $mdFileName)
$wsZhCn.Columns.Item(3).ColumnWidth = $wideColWidth
$wsZhCn.Columns.Item(9).ColumnWidth = $fullColWidth
$wsZhCn.Columns.Item(10).ColumnWidth = $fullColWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value2 = $statusText
$wsDeDe.Range("I2").Value2 = $mdFileName
$wsDeDe.Range("J2").Value2 = "ed2d9656-aa33-4b30-a561-82adbcd92828.89e1e0ed422e8a73fb20768186dace8c5457326f.de-de.xlf"
$wsDeDe.Range("K2").Value2 = "2016-08-25 19:01:41"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdUrl, "", "", $mdFileName)
$wsDeDe.Columns.Item(3).ColumnWidth = $wideColWidth
$wsDeDe.Columns.Item(9).ColumnWidth = $fullColWidth
$wsDeDe.Columns.Item(10).ColumnWidth = $fullColWidth
